# Insert 4 new weekly price rows for "Provincia de Talca" (Choclo) ahead of the
# existing block (old rows 478-513 shift down to 482-517), mirroring the
# author's commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 478:513 down by 4 rows (new block will land at 482:517).
$ws.Rows("478:481").Insert()

# Columns that stay constant for every data row in this report.
$mercadoId   = 3
$mercado     = "Femacal de La Calera"
$region      = "Coquimbo"
$codreg      = 5
$categoriaId = 100112024
$categoria   = "Choclo"
$clasif      = "Hortaliza"

# New rows: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm,
#           Unidad, Origen, Precio$/Kg, KgOUnidades
$newRows = @(
    @{ Row=478; Fecha=44578; Variedad="Choclero";           Calidad="Primera"; Vol=10000; Min=270; Max=300; Prom=286; Unidad="`$/unidad"; Origen="Provincia de Talca"; PKg=286; Q=1 },
    @{ Row=479; Fecha=44578; Variedad="Choclero";           Calidad="Segunda"; Vol=4500;  Min=200; Max=200; Prom=200; Unidad="`$/unidad"; Origen="Provincia de Talca"; PKg=200; Q=1 },
    @{ Row=480; Fecha=44578; Variedad="Dulce o Americano";  Calidad="Primera"; Vol=14500; Min=200; Max=230; Prom=214; Unidad="`$/unidad"; Origen="Provincia de Talca"; PKg=214; Q=1 },
    @{ Row=481; Fecha=44578; Variedad="Dulce o Americano";  Calidad="Segunda"; Vol=6500;  Min=150; Max=150; Prom=150; Unidad="`$/unidad"; Origen="Provincia de Talca"; PKg=150; Q=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.Fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $categoriaId
    $ws.Cells.Item($row, 7).Value2  = $categoria
    $ws.Cells.Item($row, 8).Value2  = $r.Variedad
    $ws.Cells.Item($row, 9).Value2  = $r.Calidad
    $ws.Cells.Item($row, 10).Value2 = $r.Vol
    $ws.Cells.Item($row, 11).Value2 = $r.Min
    $ws.Cells.Item($row, 12).Value2 = $r.Max
    $ws.Cells.Item($row, 13).Value2 = $r.Prom
    $ws.Cells.Item($row, 14).Value2 = $r.Unidad
    $ws.Cells.Item($row, 15).Value2 = $r.Origen
    $ws.Cells.Item($row, 16).Value2 = $r.PKg
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $clasif
}
